# Update outcome for AHC: refresh forecast mean/bounds and diff columns
# for HBV epidemic appendix data (rows 2-40, columns B,C,D,E,F,I).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 116697.836318996
$ws.Range("C2").Value = 108796.220655185
$ws.Range("D2").Value = 104952.645916324
$ws.Range("E2").Value = 124919.968583268
$ws.Range("F2").Value = 130094.722247052
$ws.Range("I2").Value = 25671.8363189964

# Row 3
$ws.Range("B3").Value = 104326.004799546
$ws.Range("C3").Value = 96508.3180848677
$ws.Range("D3").Value = 92700.1827766689
$ws.Range("E3").Value = 111770.011625949
$ws.Range("F3").Value = 115278.504384519
$ws.Range("I3").Value = 52820.0047995462

# Row 4
$ws.Range("B4").Value = 121892.686041899
$ws.Range("C4").Value = 113213.630385227
$ws.Range("D4").Value = 109819.423042992
$ws.Range("E4").Value = 130002.689920775
$ws.Range("F4").Value = 136630.010868609
$ws.Range("I4").Value = 33742.6860418994

# Row 5
$ws.Range("B5").Value = 113318.49243099
$ws.Range("C5").Value = 104200.629931617
$ws.Range("D5").Value = 100917.730737264
$ws.Range("E5").Value = 122409.620340726
$ws.Range("F5").Value = 128164.951097735
$ws.Range("I5").Value = 12056.4924309904

# Row 6
$ws.Range("B6").Value = 114519.613733229
$ws.Range("C6").Value = 105318.734066067
$ws.Range("D6").Value = 99690.4308913987
$ws.Range("E6").Value = 123622.150778757
$ws.Range("F6").Value = 129412.411617996
$ws.Range("I6").Value = 16868.6137332294

# Row 7
$ws.Range("B7").Value = 108807.511491863
$ws.Range("C7").Value = 98706.0243845419
$ws.Range("D7").Value = 94280.4284108927
$ws.Range("E7").Value = 119702.294766388
$ws.Range("F7").Value = 123847.416159765
$ws.Range("I7").Value = 9488.51149186287

# Row 8
$ws.Range("B8").Value = 112820.234006946
$ws.Range("C8").Value = 102381.771028427
$ws.Range("D8").Value = 97181.8209706125
$ws.Range("E8").Value = 123232.84348403
$ws.Range("F8").Value = 129425.431638308
$ws.Range("I8").Value = 10516.2340069457

# Row 9
$ws.Range("B9").Value = 111796.469687237
$ws.Range("C9").Value = 100899.219544795
$ws.Range("D9").Value = 94814.567904798
$ws.Range("E9").Value = 123308.356534172
$ws.Range("F9").Value = 128963.449580624
$ws.Range("I9").Value = 5661.46968723697

# Row 10
$ws.Range("B10").Value = 102629.595675564
$ws.Range("C10").Value = 90829.0915847571
$ws.Range("D10").Value = 84147.8919826373
$ws.Range("E10").Value = 114532.989332052
$ws.Range("F10").Value = 119727.863015632
$ws.Range("I10").Value = -2747.40432443614

# Row 11
$ws.Range("B11").Value = 104726.417135366
$ws.Range("C11").Value = 93344.6021751029
$ws.Range("D11").Value = 85959.0271756483
$ws.Range("E11").Value = 116593.183202054
$ws.Range("F11").Value = 122237.72374345
$ws.Range("I11").Value = 9093.41713536596

# Row 12
$ws.Range("B12").Value = 104876.237342271
$ws.Range("C12").Value = 92776.2814172262
$ws.Range("D12").Value = 86979.005158413
$ws.Range("E12").Value = 116928.197771743
$ws.Range("F12").Value = 123748.762973706
$ws.Range("I12").Value = 4315.2373422708

# Row 13
$ws.Range("B13").Value = 101613.781521199
$ws.Range("C13").Value = 87684.8087646816
$ws.Range("D13").Value = 81507.5858539772
$ws.Range("E13").Value = 113720.470208704
$ws.Range("F13").Value = 120757.047785379
$ws.Range("I13").Value = 1404.78152119939

# Row 14
$ws.Range("B14").Value = 119149.94735983
$ws.Range("C14").Value = 106569.601649526
$ws.Range("D14").Value = 97445.1430978986
$ws.Range("E14").Value = 132961.936386009
$ws.Range("F14").Value = 139178.449652833
$ws.Range("I14").Value = 20892.9473598303

# Row 15
$ws.Range("B15").Value = 107055.7128206
$ws.Range("C15").Value = 93662.067273364
$ws.Range("D15").Value = 84783.8266574478
$ws.Range("E15").Value = 121028.893537346
$ws.Range("F15").Value = 128239.31015308
$ws.Range("I15").Value = 24041.7128206

# Row 16
$ws.Range("B16").Value = 124681.748074223
$ws.Range("C16").Value = 111434.700824813
$ws.Range("D16").Value = 102224.188180103
$ws.Range("E16").Value = 138266.908391548
$ws.Range("F16").Value = 148311.042280591
$ws.Range("I16").Value = 5258.74807422332

# Row 17
$ws.Range("B17").Value = 115305.66159837
$ws.Range("C17").Value = 100801.913630374
$ws.Range("D17").Value = 92316.1814203385
$ws.Range("E17").Value = 129029.458341006
$ws.Range("F17").Value = 136380.148982306
$ws.Range("I17").Value = 4920.66159836983

# Row 18
$ws.Range("B18").Value = 116597.144031518
$ws.Range("C18").Value = 102409.719998351
$ws.Range("D18").Value = 92509.903886299
$ws.Range("E18").Value = 130988.284993117
$ws.Range("F18").Value = 140003.591276977
$ws.Range("I18").Value = 11204.1440315181

# Row 19
$ws.Range("B19").Value = 111376.005495482
$ws.Range("C19").Value = 96929.2098495315
$ws.Range("D19").Value = 88070.6294673093
$ws.Range("E19").Value = 125363.0168913
$ws.Range("F19").Value = 134264.17846838
$ws.Range("I19").Value = 9743.00549548186

# Row 20
$ws.Range("B20").Value = 115056.29437198
$ws.Range("C20").Value = 99826.5427249482
$ws.Range("D20").Value = 90872.3711421642
$ws.Range("E20").Value = 129421.550281991
$ws.Range("F20").Value = 139949.28274083
$ws.Range("I20").Value = 2769.2943719796

# Row 21
$ws.Range("B21").Value = 114683.277180079
$ws.Range("C21").Value = 99871.3285249414
$ws.Range("D21").Value = 90327.9621770071
$ws.Range("E21").Value = 130898.494640403
$ws.Range("F21").Value = 138855.521494938
$ws.Range("I21").Value = 12622.277180079

# Row 22
$ws.Range("B22").Value = 104874.759405554
$ws.Range("C22").Value = 89102.3710328425
$ws.Range("D22").Value = 82116.9868991506
$ws.Range("E22").Value = 120268.565442404
$ws.Range("F22").Value = 129628.359041328
$ws.Range("I22").Value = 3173.75940555363

# Row 23
$ws.Range("B23").Value = 106887.541214457
$ws.Range("C23").Value = 90847.5713890427
$ws.Range("D23").Value = 79931.6242424285
$ws.Range("E23").Value = 122926.673749609
$ws.Range("F23").Value = 133463.545350233
$ws.Range("I23").Value = 10945.5412144566

# Row 24
$ws.Range("B24").Value = 106396.856783558
$ws.Range("C24").Value = 89070.4501339401
$ws.Range("D24").Value = 81266.6076600328
$ws.Range("E24").Value = 122866.727630369
$ws.Range("F24").Value = 133885.911772608
$ws.Range("I24").Value = 3798.85678355781

# Row 25
$ws.Range("B25").Value = 104357.153954269
$ws.Range("C25").Value = 87327.6580561718
$ws.Range("D25").Value = 79673.5851914893
$ws.Range("E25").Value = 120543.033730009
$ws.Range("F25").Value = 131109.892776644
$ws.Range("I25").Value = -3799.84604573103

# Row 26
$ws.Range("B26").Value = 120983.653013906
$ws.Range("C26").Value = 103985.873459909
$ws.Range("D26").Value = 95129.0851412385
$ws.Range("E26").Value = 137728.133783487
$ws.Range("F26").Value = 148864.937033519
$ws.Range("I26").Value = 11452.6530139061

# Row 27
$ws.Range("B27").Value = 108906.190923625
$ws.Range("C27").Value = 91904.2142756739
$ws.Range("D27").Value = 80806.1686044425
$ws.Range("E27").Value = 126494.79352338
$ws.Range("F27").Value = 139631.800847185
$ws.Range("I27").Value = 12119.1909236248

# Row 28
$ws.Range("B28").Value = 126253.420625637
$ws.Range("C28").Value = 107384.281984111
$ws.Range("D28").Value = 98023.4117265828
$ws.Range("E28").Value = 145418.468551656
$ws.Range("F28").Value = 158072.012247063
$ws.Range("I28").Value = 12250.4206256368

# Row 29
$ws.Range("B29").Value = 117746.212815661
$ws.Range("C29").Value = 100674.314315215
$ws.Range("D29").Value = 88244.7393263178
$ws.Range("E29").Value = 137426.272495422
$ws.Range("F29").Value = 150247.179023338
$ws.Range("I29").Value = 21813.212815661

# Row 30
$ws.Range("B30").Value = 118754.849303896
$ws.Range("C30").Value = 99738.9381061177
$ws.Range("D30").Value = 88110.8669151287
$ws.Range("E30").Value = 137108.709169227
$ws.Range("F30").Value = 150942.596709944
$ws.Range("I30").Value = 15842.8493038963

# Row 31
$ws.Range("B31").Value = 113251.43580744
$ws.Range("C31").Value = 93935.708670206
$ws.Range("D31").Value = 82926.0004639524
$ws.Range("E31").Value = 133506.669933597
$ws.Range("F31").Value = 145347.812700164
$ws.Range("I31").Value = 6405.43580744034

# Row 32
$ws.Range("B32").Value = 117088.893225078
$ws.Range("C32").Value = 97667.1733422896
$ws.Range("D32").Value = 85738.2979642154
$ws.Range("E32").Value = 137411.320652011
$ws.Range("F32").Value = 148217.374007646
$ws.Range("I32").Value = 4440.8932250776

# Row 33
$ws.Range("B33").Value = 116233.370772603
$ws.Range("C33").Value = 95693.0873135125
$ws.Range("D33").Value = 85863.2461135086
$ws.Range("E33").Value = 136997.189149427
$ws.Range("F33").Value = 148769.91203827
$ws.Range("I33").Value = 858.370772603201

# Row 34
$ws.Range("B34").Value = 107138.451139737
$ws.Range("C34").Value = 87454.8366561574
$ws.Range("D34").Value = 76419.6964095279
$ws.Range("E34").Value = 128862.777791741
$ws.Range("F34").Value = 141173.572305628
$ws.Range("I34").Value = 6055.4511397373

# Row 35
$ws.Range("B35").Value = 109152.362658294
$ws.Range("C35").Value = 89993.3670640974
$ws.Range("D35").Value = 74768.2229087843
$ws.Range("E35").Value = 129318.063161517
$ws.Range("F35").Value = 145295.007726826
$ws.Range("I35").Value = 19137.362658294

# Row 36
$ws.Range("B36").Value = 108806.791514216
$ws.Range("C36").Value = 88449.4113399106
$ws.Range("D36").Value = 76313.9560216974
$ws.Range("E36").Value = 129414.062676238
$ws.Range("F36").Value = 145393.911730181
$ws.Range("I36").Value = 22435.7915142156

# Row 37
$ws.Range("B37").Value = 106331.534552175
$ws.Range("C37").Value = 84828.542531114
$ws.Range("D37").Value = 74227.4576222121
$ws.Range("E37").Value = 127780.618019491
$ws.Range("F37").Value = 142810.723260436
$ws.Range("I37").Value = 46833.5345521749

# Row 38
$ws.Range("B38").Value = 123744.725060085
$ws.Range("C38").Value = 101780.726380877
$ws.Range("D38").Value = 90879.4176298359
$ws.Range("E38").Value = 145636.297371016
$ws.Range("F38").Value = 161392.181352384
$ws.Range("I38").Value = 48954.7250600848

# Row 39
$ws.Range("B39").Value = 111463.877262871
$ws.Range("C39").Value = 89710.5890025475
$ws.Range("D39").Value = 75014.8237884718
$ws.Range("E39").Value = 134042.621974127
$ws.Range("F39").Value = 147194.431080021
$ws.Range("I39").Value = -4599.12273712868

# Row 40
$ws.Range("B40").Value = 129110.688332621
$ws.Range("C40").Value = 106367.268158415
$ws.Range("D40").Value = 92240.2154106028
$ws.Range("E40").Value = 151761.769578238
$ws.Range("F40").Value = 165180.648412675
$ws.Range("I40").Value = 2178.68833262092
